$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test case blocks for TC2 (rows 14-18), TC3 (rows 21-25) and TC4 (rows 28-32)
# had their "Steps" / "Expected Results" content reordered:
#   TC2 block now holds the "cancelar diária" content (previously in TC3 block)
#   TC3 block now holds the "detalhar diária" content (previously in TC4 block)
#   TC4 block now holds the "analisar prestação de contas" content (previously in TC2 block)

$ws.Range("B18").Value = "Beneficiário Clica em cancelar diária."
$ws.Range("D18").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"

$ws.Range("B25").Value = "Beneficiário Clica em detalhar diária."
$ws.Range("D25").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"

$ws.Range("B32").Value = "Beneficiário Clica em analisar prestação de contas."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Analisar Prestação de Contas"
